# "Made changes to the username and pwd"
#
# Adds two new columns (U: "Username", V: "pwd") to the datasheet: a
# header in row 1 and one row of sample credentials in row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1) -----------------------------------------
$ws.Range("U1").Value = "Username"
$ws.Range("V1").Value = "pwd"

# --- New data cells (row 2) --------------------------------------------
$ws.Range("U2").Value = "meghana.rao@servicemax.com"
$ws.Range("V2").Value = "cloud_101"

# --- Styling -------------------------------------------------------------
# U1 / V1 / V2 use the sheet's normal wrap-text cell style, same as every
# other plain header/data cell (e.g. A1). U2 reuses the "code" style
# already used on T1 for this kind of value. Copy formats across (instead
# of re-declaring fonts/alignment) so no new style records get minted.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("U1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("V1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("V2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("T1").Copy() | Out-Null
$ws.Range("U2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = $false

# --- Column widths for the two new columns ------------------------------
$ws.Columns.Item(21).ColumnWidth = 52.0   # column U
$ws.Columns.Item(22).ColumnWidth = 29.8   # column V

# --- View / selection ----------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 20
$ws.Range("V3").Select()
